$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.233.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.583.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  -1.37%  '
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.55'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0845'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.805.88'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.07'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.572.97'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.55'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.238.26'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '207.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('E23').Value = '  -3.32%  '
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.22'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0504'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.277.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.610'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0167'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.56%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.818'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.56'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.90%  '
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.719.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₇0983'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.38%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.23%  '
